# Set the proofing language of the (currently empty) first paragraph to
# English (United States). Word stores this as a paragraph-mark run
# property (w:pPr/w:rPr/w:lang) when the language is applied to a range
# that covers only the paragraph mark of an empty paragraph.
$d = $word.ActiveDocument
$d.Content.LanguageID = "en-US"
